$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new column before column A, shifting everything right.
$ws.Columns("A").Insert(1)

# Header cell for the new "Match ID" column
$ws.Range("A2").Value = "Match ID"

# Data rows: Match ID = 22 for every player row (4-19) and the summary row (20)
$ws.Range("A4:A19").Value = 22

# Row 20 is hidden; temporarily unhide so the write doesn't disturb its row
# height metadata, then restore the hidden state.
$ws.Rows(20).Hidden = $false
$ws.Range("A20").Value = 22
$ws.Rows(20).Hidden = $true

# Apply the bold, borderless style used for the rest of the Match ID column
$ws.Range("A2:A19").Font.Bold = $true

# Update the selection like in the authored workbook
$ws.Range("A2:A19").Select()
